$wb = $excel.ActiveWorkbook

# --- Select the Repairer sheet as active (moves tabSelected + activeTab) ---
$wsRepairer = $wb.Worksheets.Item("Repairer")
$wsRepairer.Activate()

$excel.ActiveWindow.TabRatio = 719
